$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2019, 174999, 1, "[2894]", 2894, 3977, 0.7277)
    ,@(2023, 124999, 1, "[3049]", 3049, 4058, 0.7514)
    ,@(2027, 74999, 2, "[2928, 43]", 2971, 3996, 0.7435)
    ,@(2031, 199999, 1, "[2816]", 2816, 3930, 0.7165)
    ,@(2035, 174999, 1, "[3053]", 3053, 4087, 0.747)
    ,@(2039, 124999, 1, "[2918]", 2918, 3947, 0.7393)
    ,@(2043, 174999, 1, "[2842]", 2842, 3915, 0.7259)
    ,@(2047, 99999, 1, "[3047]", 3047, 3997, 0.7623)
    ,@(2051, 174999, 1, "[3001]", 3001, 3990, 0.7521)
    ,@(2055, 99999, 1, "[2935]", 2935, 3970, 0.7393)
    ,@(2059, 199999, 1, "[2848]", 2848, 3974, 0.7167)
    ,@(2063, 74999, 1, "[2864]", 2864, 3935, 0.7278)
    ,@(2067, 99999, 1, "[2842]", 2842, 3873, 0.7338)
    ,@(2071, 324999, 1, "[2838]", 2838, 3963, 0.7161)
    ,@(2075, 74999, 1, "[2947]", 2947, 3988, 0.739)
    ,@(2079, 149999, 1, "[3010]", 3010, 4037, 0.7456)
    ,@(2083, 299999, 2, "[2909, 28]", 2937, 3949, 0.7437)
    ,@(2087, 124999, 1, "[2860]", 2860, 3990, 0.7168)
    ,@(2091, 124999, 1, "[3070]", 3070, 4056, 0.7569)
    ,@(2095, 99999, 1, "[3008]", 3008, 3977, 0.7563)
    ,@(2099, 174999, 1, "[3044]", 3044, 4057, 0.7503)
    ,@(2103, 174999, 1, "[3024]", 3024, 4042, 0.7481)
    ,@(2108, 49999, 1, "[3116]", 3116, 4104, 0.7593)
    ,@(2113, 74999, 2, "[2876, 28]", 2904, 3964, 0.7326)
    ,@(2117, 99999, 1, "[2994]", 2994, 4024, 0.744)
    ,@(2121, 74999, 1, "[2956]", 2956, 4045, 0.7308)
    ,@(2125, 149999, 1, "[2866]", 2866, 3945, 0.7265)
    ,@(2129, 174999, 1, "[2841]", 2841, 3936, 0.7218)
    ,@(2133, 99999, 1, "[2856]", 2856, 3903, 0.7317)
    ,@(2137, 99999, 2, "[2933, 31]", 2964, 4008, 0.7395)
    ,@(2141, 99999, 1, "[2886]", 2886, 3957, 0.7292999999999999)
)

$startRow = 109
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}